$wb = $excel.ActiveWorkbook

$uuid = "4fb50209-75aa-4685-83c9-d3e6b0d98908"
$mdName = "$uuid.md"
$zhXlf = "$uuid.0321b45aaf92fc8dae4bdddca1818b1d495e88b8.zh-cn.xlf"
$deXlf = "$uuid.0321b45aaf92fc8dae4bdddca1818b1d495e88b8.de-de.xlf"

$zhHandoffDt = "2016-03-08 06:25:19"
$zhHandbackDt = "2016-03-08 06:26:05"
$deHandoffDt = "2016-03-08 06:25:28"
$deHandbackDt = "2016-03-08 06:26:21"

$status = "Handed back: in sync with en-US"
$reason = "Include"

# ---------------------------------------------------------------------------
# Sheet "Overview" : add row 4 (File Name | zh-cn | de-de)
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

# Duplicate the formatting of the previous data row (row 3) by copying it
# down - this keeps the existing cell styles (hyperlink style, date style,
# ...) intact instead of fabricating brand new ones.
$wsOverview.Rows(3).Copy()
$wsOverview.Rows(4).Insert()

$wsOverview.Range("A4").Value = $mdName
$wsOverview.Range("B4").Value = $status
$wsOverview.Range("C4").Value = $status

$wsOverview.Hyperlinks.Add(
    $wsOverview.Range("A4"),
    "https://github.com/OpenLocalizationTest/oltest/blob/0321b45aaf92fc8dae4bdddca1818b1d495e88b8/e2e/$mdName",
    "",
    "",
    $mdName
) | Out-Null

# ---------------------------------------------------------------------------
# Sheet "zh-cn" : add row 4
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Rows(3).Copy()
$wsZh.Rows(4).Insert()

$wsZh.Range("A4").Value = $mdName
$wsZh.Range("B4").Value = $status
$wsZh.Range("C4").Value = $zhXlf
$wsZh.Range("D4").Value = $zhHandoffDt
$wsZh.Range("E4").Value = $mdName
$wsZh.Range("F4").Value = $zhXlf
$wsZh.Range("G4").Value = $zhHandbackDt
$wsZh.Range("H4").Value = $reason

$wsZh.Hyperlinks.Add(
    $wsZh.Range("A4"),
    "https://github.com/OpenLocalizationTest/oltest/blob/0321b45aaf92fc8dae4bdddca1818b1d495e88b8/e2e/$mdName",
    "",
    "",
    $mdName
) | Out-Null

$wsZh.Hyperlinks.Add(
    $wsZh.Range("C4"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0321b45aaf92fc8dae4bdddca1818b1d495e88b8/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/$zhXlf",
    "",
    "",
    $zhXlf
) | Out-Null

$wsZh.Hyperlinks.Add(
    $wsZh.Range("E4"),
    "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/0321b45aaf92fc8dae4bdddca1818b1d495e88b8/e2e/$mdName",
    "",
    "",
    $mdName
) | Out-Null

$wsZh.Hyperlinks.Add(
    $wsZh.Range("F4"),
    "https://github.com/OpenLocalizationTestOrg/olhandback/blob/0321b45aaf92fc8dae4bdddca1818b1d495e88b8/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/$zhXlf",
    "",
    "",
    $zhXlf
) | Out-Null

# ---------------------------------------------------------------------------
# Sheet "de-de" : add row 4
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Rows(3).Copy()
$wsDe.Rows(4).Insert()

$wsDe.Range("A4").Value = $mdName
$wsDe.Range("B4").Value = $status
$wsDe.Range("C4").Value = $deXlf
$wsDe.Range("D4").Value = $deHandoffDt
$wsDe.Range("E4").Value = $mdName
$wsDe.Range("F4").Value = $deXlf
$wsDe.Range("G4").Value = $deHandbackDt
$wsDe.Range("H4").Value = $reason

$wsDe.Hyperlinks.Add(
    $wsDe.Range("A4"),
    "https://github.com/OpenLocalizationTest/oltest/blob/0321b45aaf92fc8dae4bdddca1818b1d495e88b8/e2e/$mdName",
    "",
    "",
    $mdName
) | Out-Null

$wsDe.Hyperlinks.Add(
    $wsDe.Range("C4"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0321b45aaf92fc8dae4bdddca1818b1d495e88b8/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/$deXlf",
    "",
    "",
    $deXlf
) | Out-Null

$wsDe.Hyperlinks.Add(
    $wsDe.Range("E4"),
    "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/0321b45aaf92fc8dae4bdddca1818b1d495e88b8/e2e/$mdName",
    "",
    "",
    $mdName
) | Out-Null

$wsDe.Hyperlinks.Add(
    $wsDe.Range("F4"),
    "https://github.com/OpenLocalizationTestOrg/olhandback/blob/0321b45aaf92fc8dae4bdddca1818b1d495e88b8/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/$deXlf",
    "",
    "",
    $deXlf
) | Out-Null
